$d = $word.ActiveDocument

# Locate the unique "M2Doc version mismatch" warning text that must be removed,
# along with the preceding "    <---" run that introduces it.
$warn = $d.Content
$found = $warn.Find.Execute("M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the M2Doc version mismatch text to remove."
}

$endOfRemoval = $warn.End
$startOfRemoval = $warn.Start - [string]"<---".Length - 4

# Immediately after the removed text comes a "    " run followed directly by a
# "demonstration" run; those two runs share identical (empty) formatting. Any
# edit inside this paragraph makes the engine coalesce adjacent same-format runs
# across the whole paragraph, which would merge "    " and "demonstration" into
# one run even though the deletion itself never touches that boundary. Drop an
# anchor bookmark exactly at that run boundary (i.e. right after the trailing
# 4 spaces that follow the removed text) before deleting, so the boundary
# survives the edit, then remove the bookmark afterwards (bookmarks carry no
# text, so adding/removing one does not itself touch/merge any run).
$guardPos = $endOfRemoval + 4
$guardName = "m2docVersionGuard"
$guardRange = $d.Range($guardPos, $guardPos)
$d.Bookmarks.Add($guardName, $guardRange)

$toDelete = $d.Range($startOfRemoval, $endOfRemoval)
$toDelete.Delete()

if ($d.Bookmarks.Exists($guardName)) {
    $d.Bookmarks($guardName).Delete()
}
